$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.934.58"
$ws.Range("E2").Value = "  +1.67%  "

$ws.Range("D3").Value = "2.641.95"
$ws.Range("E3").Value = "  +0.73%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.35"
$ws.Range("E5").Value = "  -0.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.40"
$ws.Range("E6").Value = "  +2.22%  "

$ws.Range("E7").Value = "  -0.16%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.591"
$ws.Range("E8").Value = "  +0.30%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.115"
$ws.Range("E9").Value = "  +5.29%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.398"
$ws.Range("E10").Value = "  +2.84%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.80"
$ws.Range("E11").Value = "  +1.58%  "

$ws.Range("E12").Value = "  +1.23%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.93"
$ws.Range("E13").Value = "  +3.15%  "

$ws.Range("D14").Value = "3.105.07"
$ws.Range("E14").Value = "  +0.30%  "

$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "64.683.28"
$ws.Range("E15").Value = "  +1.52%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000171"
$ws.Range("E16").Value = "  +10.61%  "

$ws.Range("D17").Value = "2.604.38"
$ws.Range("E17").Value = "  -1.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.51"
$ws.Range("E18").Value = "  +0.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.81"
$ws.Range("E19").Value = "  +1.64%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "352.40"
$ws.Range("E20").Value = "  +0.95%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.23"
$ws.Range("E21").Value = "  +4.43%  "

$ws.Range("E22").Value = "  +0.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.85"
$ws.Range("E23").Value = "  +0.61%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.73"
$ws.Range("E24").Value = "  +0.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.45"
$ws.Range("E25").Value = "  +0.16%  "

$ws.Range("E26").Value = "  -3.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.35"
$ws.Range("E27").Value = "  +3.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.164"
$ws.Range("E28").Value = "  +0.95%  "

$ws.Range("E29").Value = "  +0.19%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "531.51"
$ws.Range("E30").Value = "  -5.10%  "

$ws.Range("D31").Value = "0.0₃0917"
$ws.Range("E31").Value = "  +7.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.05"
$ws.Range("E32").Value = "  -0.48%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.82"
$ws.Range("E33").Value = "  +3.17%  "

$ws.Range("E34").Value = "  +8.40%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.37"
$ws.Range("E35").Value = "  +1.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.427"
$ws.Range("E36").Value = "  +2.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "165.04"
$ws.Range("E37").Value = "  -1.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.03"
$ws.Range("E38").Value = "  +3.48%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "20.20"
$ws.Range("E39").Value = "  +2.73%  "

$ws.Range("E40").Value = "  -0.10%  "

$ws.Range("E41").Value = "  -0.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "167.31"
$ws.Range("E42").Value = "  -0.31%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.63"
$ws.Range("E43").Value = "  +4.49%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.10"
$ws.Range("E44").Value = "  +2.79%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0605"
$ws.Range("E45").Value = "  +2.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.27"
$ws.Range("E46").Value = "  +5.27%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.24"
$ws.Range("E47").Value = "  +8.66%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.643"
$ws.Range("E48").Value = "  +1.09%  "

$ws.Range("E49").Value = "  -0.77%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0983"
$ws.Range("E50").Value = "  +1.32%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.43"
$ws.Range("E51").Value = "  -0.01%  "
